$wb = $excel.ActiveWorkbook

# --- 1. Update DOI version referenced on the "glossary" sheet ---
$glossary = $wb.Worksheets.Item("glossary")
$glossary.Range("A2").Value = "DOI: 10.25573/serc.14555511.v3"

# --- 2. Broaden "sheet name" for the sample_collection_day/month/year rows ---
#        (rows 16-18) so they also apply to "transect metadata", not just
#        "site metadata".
$glossary.Range("F16").Value = "site metadata, transect metadata"
$glossary.Range("F17").Value = "site metadata, transect metadata"
$glossary.Range("F18").Value = "site metadata, transect metadata"

# --- 3. Add sample_collection_year / _month / _day columns to the
#        "transect metadata" sheet, ahead of the transect begin/end
#        lat-long columns. ---
$transect = $wb.Worksheets.Item("transect metadata")

# Insert three new blank columns at D (pushes the existing D:J -> G:M,
# carrying their formatting/widths along for the ride).
$transect.Range("D1:F1").EntireColumn.Insert()

$transect.Range("D1").Value = "sample_collection_year"
$transect.Range("E1").Value = "sample_collection_month"
$transect.Range("F1").Value = "sample_collection_day"

# Size the three new columns.
$transect.Columns.Item(4).ColumnWidth = 21.8333
$transect.Columns.Item(5).ColumnWidth = 22.8333
$transect.Columns.Item(6).ColumnWidth = 20.8333
